$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = "232.89"
    "D3"  = "23.05"
    "D4"  = "5.551"
    "D5"  = "0.05638"
    "D6"  = "3.416"
    "D7"  = "6.482"
    "D8"  = "1.271"
    "D9"  = "0.8011"
    "D10" = "0.1421"
    "D11" = "0.07505"
    "D12" = "0.03247"
    "D13" = "0.02930"
    "D14" = "0.09226"
    "D15" = "0.001678"
    "D16" = "3.269"
    "D17" = "0.04730"
    "D18" = "0.0005988"
    "D19" = "0.006230"
    "D20" = "0.005355"
    "D21" = "0.001068"
    "D23" = "3.689"
    "D26" = "0.1279"
    "D27" = "0.0006747"
    "D40" = "0.04134"
    "D41" = "0.007044"
    "D42" = "0.003455"
    "D43" = "0.1048"
    "D44" = "0.008340"
    "D46" = "0.00005575"
    "D47" = "0.00000000752"
    "D48" = "0.7875"
    "D49" = "0.09644"
    "D50" = "0.00002107"
    "D51" = "0.01013"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
